# Auto-generated edit script: numeric cell updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# Mirrors the OOXML diff: currentAveragePrice / NQ / HQ / LevePrice / LeveProfit recalculated values,
# plus a handful of previously-empty LeveProfitHQ (N) / LeveProfitNQ (M) cells that are now populated.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1246.3889
$ws.Range("I132").Value = 1060
$ws.Range("K132").Value = 3180
$ws.Range("M132").Value = -650
$ws.Range("H137").Value = 999.1579
$ws.Range("I137").Value = 940.4706
$ws.Range("J137").Value = 1498
$ws.Range("K137").Value = 2821.4118
$ws.Range("L137").Value = 4494
$ws.Range("M137").Value = -271.4117999999999
$ws.Range("N137").Value = -9594
$ws.Range("H138").Value = 4206.4707
$ws.Range("I138").Value = 7546.3335
$ws.Range("J138").Value = 3490.7856
$ws.Range("K138").Value = 22639.0005
$ws.Range("L138").Value = 10472.3568
$ws.Range("M138").Value = -17499.0005
$ws.Range("N138").Value = -20752.3568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3357.4517
$ws.Range("I2").Value = 2503.9524
$ws.Range("K2").Value = 2503.9524
$ws.Range("M2").Value = -2390.9524
$ws.Range("H32").Value = 6000.686
$ws.Range("I32").Value = 5251.396
$ws.Range("K32").Value = 5251.396
$ws.Range("M32").Value = -4964.396
$ws.Range("H45").Value = 122790.35
$ws.Range("I45").Value = 225719.33
$ws.Range("K45").Value = 225719.33
$ws.Range("M45").Value = -225342.33
$ws.Range("H50").Value = 32622
$ws.Range("J50").Value = 49996
$ws.Range("L50").Value = 49996
$ws.Range("N50").Value = -51424
$ws.Range("H61").Value = 7869.6787
$ws.Range("I61").Value = 6572.95
$ws.Range("K61").Value = 6572.95
$ws.Range("M61").Value = -6360.95
$ws.Range("H86").Value = 25017500
$ws.Range("J86").Value = 35000
$ws.Range("L86").Value = 35000
$ws.Range("N86").Value = -37372
$ws.Range("H88").Value = 1345.4
$ws.Range("I88").Value = 1254.8572
$ws.Range("K88").Value = 1254.8572
$ws.Range("M88").Value = -848.8571999999999
$ws.Range("H89").Value = 25017500
$ws.Range("J89").Value = 35000
$ws.Range("L89").Value = 105000
$ws.Range("N89").Value = -116856
$ws.Range("H91").Value = 1345.4
$ws.Range("I91").Value = 1254.8572
$ws.Range("K91").Value = 1254.8572
$ws.Range("M91").Value = 149.1428000000001
$ws.Range("H116").Value = 3357.4517
$ws.Range("I116").Value = 2503.9524
$ws.Range("K116").Value = 2503.9524
$ws.Range("M116").Value = -209.9524000000001
$ws.Range("H132").Value = 4973.2812
$ws.Range("I132").Value = 4239.9614
$ws.Range("J132").Value = 8151
$ws.Range("K132").Value = 12719.8842
$ws.Range("L132").Value = 24453
$ws.Range("M132").Value = -10189.8842
$ws.Range("N132").Value = -29513
$ws.Range("H136").Value = 7869.6787
$ws.Range("I136").Value = 6572.95
$ws.Range("K136").Value = 19718.85
$ws.Range("M136").Value = -17168.85

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3357.4517
$ws.Range("I3").Value = 2503.9524
$ws.Range("K3").Value = 2503.9524
$ws.Range("M3").Value = -2389.9524
$ws.Range("H20").Value = 3445.25
$ws.Range("I20").Value = 3819.25
$ws.Range("J20").Value = 2697.25
$ws.Range("K20").Value = 3819.25
$ws.Range("L20").Value = 2697.25
$ws.Range("M20").Value = -3572.25
$ws.Range("N20").Value = -3191.25
$ws.Range("H86").Value = 33339916
$ws.Range("I86").Value = 8875.5
$ws.Range("J86").Value = 100002000
$ws.Range("K86").Value = 8875.5
$ws.Range("L86").Value = 100002000
$ws.Range("M86").Value = -7752.5
$ws.Range("N86").Value = -100004246
$ws.Range("H89").Value = 33339916
$ws.Range("I89").Value = 8875.5
$ws.Range("J89").Value = 100002000
$ws.Range("K89").Value = 44377.5
$ws.Range("L89").Value = 500010000
$ws.Range("M89").Value = -38761.5
$ws.Range("N89").Value = -500021232
$ws.Range("H134").Value = 6963.3794
$ws.Range("I134").Value = 6909.346
$ws.Range("J134").Value = 7431.6665
$ws.Range("K134").Value = 20728.038
$ws.Range("L134").Value = 22294.9995
$ws.Range("M134").Value = -18193.038
$ws.Range("N134").Value = -27364.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 483.6
$ws.Range("I22").Value = 354.5
$ws.Range("K22").Value = 354.5
$ws.Range("M22").Value = -4.5
$ws.Range("H31").Value = 2789.6226
$ws.Range("J31").Value = 3680.2632
$ws.Range("L31").Value = 3680.2632
$ws.Range("N31").Value = -4270.263199999999
$ws.Range("H34").Value = 2789.6226
$ws.Range("J34").Value = 3680.2632
$ws.Range("L34").Value = 3680.2632
$ws.Range("N34").Value = -4084.2632
$ws.Range("H115").Value = 39290
$ws.Range("J115").Value = 39290
$ws.Range("L115").Value = 39290
$ws.Range("N115").Value = -41640
$ws.Range("H132").Value = 3368.25
$ws.Range("I132").Value = 1907.1538
$ws.Range("K132").Value = 5721.4614
$ws.Range("M132").Value = -3191.4614
$ws.Range("H134").Value = 10132.333
$ws.Range("I134").Value = 7499.25
$ws.Range("K134").Value = 22497.75
$ws.Range("M134").Value = -19962.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 141.83333
$ws.Range("I14").Value = 141.83333
$ws.Range("K14").Value = 425.49999
$ws.Range("M14").Value = -252.49999
$ws.Range("H122").Value = 2495.2727
$ws.Range("I122").Value = 2007
$ws.Range("J122").Value = 3349.75
$ws.Range("K122").Value = 18063
$ws.Range("L122").Value = 30147.75
$ws.Range("M122").Value = -15613
$ws.Range("N122").Value = -35047.75
$ws.Range("H131").Value = 29413668
$ws.Range("J131").Value = 2883.9
$ws.Range("L131").Value = 8651.700000000001
$ws.Range("N131").Value = -18731.7
$ws.Range("H132").Value = 33333918
$ws.Range("I132").Value = 1004
$ws.Range("K132").Value = 9036
$ws.Range("M132").Value = -6506
$ws.Range("H137").Value = 8219.764999999999
$ws.Range("I137").Value = 2137.4
$ws.Range("J137").Value = 16908.857
$ws.Range("K137").Value = 6412.200000000001
$ws.Range("L137").Value = 50726.571
$ws.Range("M137").Value = -1312.200000000001
$ws.Range("N137").Value = -60926.571

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4166.1665
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 4166.1665
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008
$ws.Range("H97").Value = 789.2105
$ws.Range("I97").Value = 743.5714
$ws.Range("K97").Value = 743.5714
$ws.Range("M97").Value = -247.5714
$ws.Range("H126").Value = 5751.1
$ws.Range("I126").Value = 6351.375
$ws.Range("J126").Value = 3350
$ws.Range("K126").Value = 19054.125
$ws.Range("L126").Value = 10050
$ws.Range("M126").Value = -16584.125
$ws.Range("N126").Value = -14990
$ws.Range("H132").Value = 2051.4517
$ws.Range("I132").Value = 1895.4445
$ws.Range("K132").Value = 5686.333500000001
$ws.Range("M132").Value = -3156.333500000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4164.4287
$ws.Range("I40").Value = 3725.25
$ws.Range("K40").Value = 3725.25
$ws.Range("M40").Value = -3589.25
$ws.Range("H100").Value = 3128251
$ws.Range("I100").Value = 6252613
$ws.Range("J100").Value = 3888.75
$ws.Range("K100").Value = 6252613
$ws.Range("L100").Value = 3888.75
$ws.Range("M100").Value = -6252072
$ws.Range("N100").Value = -4970.75
$ws.Range("H122").Value = 3252.5557
$ws.Range("I122").Value = 3252.5557
$ws.Range("K122").Value = 9757.667099999999
$ws.Range("M122").Value = -7307.667099999999
$ws.Range("H132").Value = 12911.714
$ws.Range("I132").Value = 14695.777
$ws.Range("K132").Value = 44087.331
$ws.Range("M132").Value = -41557.331
$ws.Range("H136").Value = 2876.6316
$ws.Range("J136").Value = 3474.8
$ws.Range("L136").Value = 10424.4
$ws.Range("N136").Value = -15524.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 1082
$ws.Range("I55").Value = 874
$ws.Range("K55").Value = 874
$ws.Range("M55").Value = -597
$ws.Range("H62").Value = 4496
$ws.Range("J62").Value = 4496
$ws.Range("L62").Value = 4496
$ws.Range("N62").Value = -5744
$ws.Range("H65").Value = 4496
$ws.Range("J65").Value = 4496
$ws.Range("L65").Value = 22480
$ws.Range("N65").Value = -28720
$ws.Range("H81").Value = 2689.5833
$ws.Range("J81").Value = 5615.2
$ws.Range("L81").Value = 11230.4
$ws.Range("N81").Value = -13352.4
$ws.Range("H84").Value = 2689.5833
$ws.Range("J84").Value = 5615.2
$ws.Range("L84").Value = 56152
$ws.Range("N84").Value = -66760
$ws.Range("H107").Value = 909.1429000000001
$ws.Range("I107").Value = 909.1429000000001
$ws.Range("K107").Value = 2727.4287
$ws.Range("M107").Value = -807.4287000000004
